# Slide 3 ("X3DJSONLD Modules") -> "Content Placeholder 2":
# add explanatory "server-side"/"client-side" notes to three of the bullets
# (commit: "Added server side and client side references").
#
# NOTE on the <a:normAutofit fontScale="92500" lnSpcReduction="10000"/> ->
# lnSpcReduction="20000" tweak in the canonical diff: that percentage is a
# value PowerPoint's live text layout engine derives purely from rendering
# metrics (how much extra line-spacing squeeze is needed once the bullet
# text grows) and there is no AutoSize/NormAutofit scale or line-spacing
# property on TextFrame/TextFrame2/ParagraphFormat (real PowerPoint doesn't
# expose one either) to set it explicitly. Touching TextFrame.AutoSize here
# only resets <a:normAutofit> to a bare, attribute-less element, which is
# further from the target than simply leaving the existing
# fontScale="92500" lnSpcReduction="10000" in place, so it is intentionally
# left untouched and only the text runs below are edited.

$dash = [char]0x2013   # "-" en dash, matches the existing bullet typography
$apos = [char]0x2019   # "'" right single quotation mark

$pres  = $ppt.ActivePresentation
$slide = $pres.Slides.Item(3)
$shape = $slide.Shapes.Item(2)          # "Content Placeholder 2"
$tr    = $shape.TextFrame.TextRange

# --- Bullet 1: "X3DJSONLD.js ... independent of Jquery." -------------------
# Runs: "X3DJSONLD.js ... independent of " | "Jquery" (err=1) | "."
# Only the trailing "." run gains text; addressed by run index so the two
# earlier runs (and the err="1" flag on "Jquery") are left exactly as-is.
$para1 = $tr.Paragraphs(1, 1)
$tailRun1 = $para1.Runs(3, 1)
$tailRun1.Text = "X"                     # break the shared-prefix diff match
$tailRun1 = $para1.Runs(3, 1)
$tailRun1.Text = ".  Server and client side.  Contains some client-side code which shouldn${apos}t be used on server."

# --- Bullet 2: "loaderJQuery.js ..." ----------------------------------------
# Single-run paragraph -> replace the whole run. The two-step (dummy value,
# then the final text) keeps the result a single run instead of PowerPoint
# splitting it into an unchanged-prefix / new-suffix run pair.
$para2 = $tr.Paragraphs(2, 1)
$para2.Text = "X"
$para2 = $tr.Paragraphs(2, 1)
$para2.Text = "loaderJQuery.js ${dash} jQuery and other integrations useful for a web page.  Client-side"

# --- Bullet 3: "convertJSON.js ... Also contains JSON validator." ----------
# Runs: "convertJSON.js ... send to " | "serializer" (err=1) | ".  Also ..."
# Same run-index targeting as bullet 1.
$para3 = $tr.Paragraphs(3, 1)
$tailRun3 = $para3.Runs(3, 1)
$tailRun3.Text = "X"
$tailRun3 = $para3.Runs(3, 1)
$tailRun3.Text = ".  Also contains JSON validator.  Server-side."
